# Applies the changes described by the commit:
#  - Refresh the cached "auto date" placeholder text on the slide master,
#    all slide layouts and the notes master (08.07.2025 -> 16.07.2025 /
#    7/8/2025 -> 7/16/2025).
#  - Slide 1: subtitle "Created by: Said Cetin" -> "Said Cetin".
#  - Slide 2: title "Agenda" -> "Kapitel 1 Agenda".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*" -or $shape.Name -like "Datumsplatzhalter*") {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $text = $tf.TextRange.Text
                if ($text -eq "7/8/2025") {
                    $tf.TextRange.Text = "7/16/2025"
                } elseif ($text -eq "08.07.2025") {
                    $tf.TextRange.Text = "16.07.2025"
                }
            }
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout attached to the master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master date placeholder.
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# Slide 1: subtitle text simplified.
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "Created by: Said Cetin") {
            $shape.TextFrame.TextRange.Text = "Said Cetin"
        }
    }
}

# Slide 2: title text updated.
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shape = $slide2.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "Agenda") {
            $shape.TextFrame.TextRange.Text = "Kapitel 1 Agenda"
        }
    }
}
